$wb = $excel.ActiveWorkbook

# --- DatosCuenta (sheet1) ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "PruebaDecima"
$wsCuenta.Range("B2").Value = "PruebaDecima"
$wsCuenta.Range("C2").Value = 270100100
$wsCuenta.Range("D2").Value = 105

# --- DatosHogar (sheet2) ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 626

# --- DatosMotor (sheet3) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA007"
$wsMotor.Range("B2").Value = "ABC12SSMA007"
$wsMotor.Range("C2").Value = "ZAZ123SSMA007"

# --- DatosAP (sheet4) ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200106

# --- Selections / active sheet order matters: set non-final sheets first ---
$wsHogar.Range("A2").Select()
$wsMotor.Range("D5").Select()
$wsCuenta.Range("D3").Select()
$wsAP.Range("A2").Select()
